# Access the active workbook/worksheet (already open per harness contract)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the "contraste"/interaction columns K..Q (shared strings 10..16)
$ws.Range("K1").Value = "1+234"
$ws.Range("L1").Value = "2+134"
$ws.Range("M1").Value = "3+124"
$ws.Range("N1").Value = "4+123"
$ws.Range("O1").Value = "12+34"
$ws.Range("P1").Value = "13+24"
$ws.Range("Q1").Value = "14+23"

# Fill in the interaction/contrast formulas for data rows 2..12
# K = A (v1), L = C (v2), M = E (v3), N = G (v4)
# O = K*L, P = K*M, Q = K*N
for ($row = 2; $row -le 12; $row++) {
    $ws.Range("K$row").Formula = "=A$row"
    $ws.Range("L$row").Formula = "=C$row"
    $ws.Range("M$row").Formula = "=E$row"
    $ws.Range("N$row").Formula = "=G$row"
    $ws.Range("O$row").Formula = "=K$row*L$row"
    $ws.Range("P$row").Formula = "=K$row*M$row"
    $ws.Range("Q$row").Formula = "=K$row*N$row"
}

# Update the saved selection to match the author's final cursor position
$ws.Range("N16").Select()
